# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for the rows
# whose dialog-act annotations changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3;   I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 17;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 19;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 21;  I = "%";  J = "Uninterpretable" },
    @{ Row = 35;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 46;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 51;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 56;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 64;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 65;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 77;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 106; I = "sv"; J = "Statement-opinion" },
    @{ Row = 122; I = "sv"; J = "Statement-opinion" },
    @{ Row = 125; I = "ba"; J = "Appreciation" },
    @{ Row = 129; I = "sv"; J = "Statement-opinion" },
    @{ Row = 149; I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 152; I = "qy"; J = "Yes-No-Question" },
    @{ Row = 154; I = "sd"; J = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
